$p = $ppt.ActivePresentation

# Slide 7's notes page body placeholder: consolidate the "This is a blank..." runs
# into a single run. Setting the text to a placeholder value first forces the
# engine to actually rebuild the paragraph (a no-op assignment of the same
# text is a short-circuited no-op), then we set it back to the real text.
$s7 = $p.Slides.Item(7)
$np = $s7.NotesPage
$notesTr = $np.Shapes.Item(2).TextFrame.TextRange
$notesTr.Text = "x"
$notesTr.Text = "This is a blank slide: does it have a footer?"

# Title placeholders on slides 2, 4, 5, 6: consolidate "Slide" " " "N" runs
# into a single run, using the same placeholder-then-real-value trick.
$s2 = $p.Slides.Item(2)
$tr2 = $s2.Shapes.Item(1).TextFrame.TextRange
$tr2.Text = "x"
$tr2.Text = "Slide 1"

$s4 = $p.Slides.Item(4)
$tr4 = $s4.Shapes.Item(1).TextFrame.TextRange
$tr4.Text = "x"
$tr4.Text = "Slide 3"

$s5 = $p.Slides.Item(5)
$tr5 = $s5.Shapes.Item(1).TextFrame.TextRange
$tr5.Text = "x"
$tr5.Text = "Slide 4"

$s6 = $p.Slides.Item(6)
$tr6 = $s6.Shapes.Item(1).TextFrame.TextRange
$tr6.Text = "x"
$tr6.Text = "Slide 5"
